$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = 1988393
$ws.Range("B2").Value = 44165

# Row 3 (A3 loses its special style -> becomes plain/Normal)
$ws.Range("A3").Value = 1999934
$ws.Range("B3").Value = 95656
$ws.Range("A3").Style = "Normal"

# Row 4 (A4 loses its special style -> becomes plain/Normal)
$ws.Range("A4").Value = 1999982
$ws.Range("B4").Value = 51980
$ws.Range("A4").Style = "Normal"

# Row 5 (newly filled in, A5 has no special style)
$ws.Range("A5").Value = 2012994
$ws.Range("B5").Value = 37976
$ws.Range("C5").Value = -7
$ws.Range("D5").Value = 0
$ws.Range("E5").Value = 3.5
$ws.Range("A5").Style = "Normal"

# Row 6 (newly filled in, A6 has no special style)
$ws.Range("A6").Value = 2045531
$ws.Range("B6").Value = 19136
$ws.Range("C6").Value = -7
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 3.5
$ws.Range("A6").Style = "Normal"

# Shrink the conditional formatting range from A2:A3 down to just A2
$fcs = $ws.Range("A2").FormatConditions
for ($i = 1; $i -le $fcs.Count; $i++) {
    $fcs.Item($i).ModifyAppliesToRange($ws.Range("A2"))
}

# Update the active selection
$ws.Range("C13").Select()
